# Capitalize the categorical rating values (high/medium/low/best/worst)
# throughout Sheet1's density_by_city, air_quality_by_city,
# unemployment_by_city, linkjobs_by_city and expenses_by_city columns
# (E, G, I, K, P) to High/Medium/Low/Best/Worst.
#
# Terms are replaced one at a time (all occurrences of "high" first, then
# "medium", then "low", then "best", then "worst") so that the workbook's
# shared-string table is rebuilt on save in that same encounter order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.UsedRange.Rows.Count

$terms = @("high", "medium", "low", "best", "worst")
$replacements = @{
    "high"   = "High"
    "medium" = "Medium"
    "low"    = "Low"
    "best"   = "Best"
    "worst"  = "Worst"
}
$columns = @("E", "G", "I", "K", "P")

foreach ($term in $terms) {
    for ($r = 2; $r -le $lastRow; $r++) {
        foreach ($col in $columns) {
            $cell = $ws.Range("$col$r")
            $cur = $cell.Value2
            if ($cur -eq $term) {
                $cell.Value = $replacements[$term]
            }
        }
    }
}
